$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 51.3
$ws.Range("I5").Value = 51.3
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 51.3
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 63.7
$ws.Range("N5").ClearContents()
$ws.Range("H19").Value = 1246.8334
$ws.Range("I19").Value = 1400
$ws.Range("J19").Value = 1232.909
$ws.Range("K19").Value = 1400
$ws.Range("L19").Value = 1232.909
$ws.Range("M19").Value = -1225
$ws.Range("N19").Value = -1582.909
$ws.Range("H32").Value = 2066.9092
$ws.Range("I32").Value = 1831
$ws.Range("K32").Value = 1831
$ws.Range("M32").Value = -1505
$ws.Range("H33").Value = 804.75
$ws.Range("I33").Value = 573
$ws.Range("J33").Value = 1500
$ws.Range("K33").Value = 573
$ws.Range("L33").Value = 1500
$ws.Range("M33").Value = -344
$ws.Range("N33").Value = -1958
$ws.Range("H80").Value = 999.5
$ws.Range("J80").Value = 999.5
$ws.Range("L80").Value = 2998.5
$ws.Range("N80").Value = -4994.5
$ws.Range("H83").Value = 999.5
$ws.Range("J83").Value = 999.5
$ws.Range("L83").Value = 8995.5
$ws.Range("N83").Value = -18979.5
$ws.Range("H88").Value = 2708.25
$ws.Range("J88").Value = 2611
$ws.Range("L88").Value = 2611
$ws.Range("N88").Value = -3423
$ws.Range("H91").Value = 2708.25
$ws.Range("J91").Value = 2611
$ws.Range("L91").Value = 2611
$ws.Range("N91").Value = -5419
$ws.Range("H106").Value = 3802995
$ws.Range("I106").Value = 4800499.5
$ws.Range("J106").Value = 12479.4
$ws.Range("K106").Value = 4800499.5
$ws.Range("L106").Value = 12479.4
$ws.Range("M106").Value = -4799868.5
$ws.Range("N106").Value = -13741.4
$ws.Range("H116").Value = 11649.342
$ws.Range("I116").Value = 14561.875
$ws.Range("J116").Value = 10943.272
$ws.Range("K116").Value = 14561.875
$ws.Range("L116").Value = 10943.272
$ws.Range("M116").Value = -11119.875
$ws.Range("N116").Value = -17827.272
$ws.Range("H127").Value = 1709.3529
$ws.Range("I127").Value = 1551.091
$ws.Range("K127").Value = 4653.272999999999
$ws.Range("M127").Value = 306.7270000000008
$ws.Range("H131").Value = 2794.0588
$ws.Range("I131").Value = 1778.2142
$ws.Range("K131").Value = 5334.642599999999
$ws.Range("M131").Value = -294.6425999999992
$ws.Range("H137").Value = 1863.2941
$ws.Range("I137").Value = 927
$ws.Range("K137").Value = 2781
$ws.Range("M137").Value = -231
$ws.Range("H138").Value = 2773.3333
$ws.Range("I138").Value = 2128.9678
$ws.Range("J138").Value = 3344.0571
$ws.Range("K138").Value = 6386.903399999999
$ws.Range("L138").Value = 10032.1713
$ws.Range("M138").Value = -1246.903399999999
$ws.Range("N138").Value = -20312.1713

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2263.5322
$ws.Range("I32").Value = 1868.4108
$ws.Range("K32").Value = 1868.4108
$ws.Range("M32").Value = -1581.4108
$ws.Range("H61").Value = 9679.526
$ws.Range("I61").Value = 6766.6924
$ws.Range("K61").Value = 6766.6924
$ws.Range("M61").Value = -6554.6924
$ws.Range("H88").Value = 8157.4375
$ws.Range("I88").Value = 2035
$ws.Range("J88").Value = 9032.071
$ws.Range("K88").Value = 2035
$ws.Range("L88").Value = 9032.071
$ws.Range("M88").Value = -1629
$ws.Range("N88").Value = -9844.071
$ws.Range("H91").Value = 8157.4375
$ws.Range("I91").Value = 2035
$ws.Range("J91").Value = 9032.071
$ws.Range("K91").Value = 2035
$ws.Range("L91").Value = 9032.071
$ws.Range("M91").Value = -631
$ws.Range("N91").Value = -11840.071
$ws.Range("H110").Value = 2063.476
$ws.Range("I110").Value = 1868.5555
$ws.Range("K110").Value = 1868.5555
$ws.Range("M110").Value = 176.4445000000001
$ws.Range("H136").Value = 9679.526
$ws.Range("I136").Value = 6766.6924
$ws.Range("K136").Value = 20300.0772
$ws.Range("M136").Value = -17750.0772

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 50236.25
$ws.Range("J82").Value = 89997.5
$ws.Range("L82").Value = 89997.5
$ws.Range("N82").Value = -90763.5
$ws.Range("H85").Value = 50236.25
$ws.Range("J85").Value = 89997.5
$ws.Range("L85").Value = 89997.5
$ws.Range("N85").Value = -92649.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 308.83334
$ws.Range("I7").Value = 304
$ws.Range("K7").Value = 304
$ws.Range("M7").Value = -191
$ws.Range("H62").Value = 13866.286
$ws.Range("I62").Value = 6350.8335
$ws.Range("J62").Value = 19502.875
$ws.Range("K62").Value = 6350.8335
$ws.Range("L62").Value = 19502.875
$ws.Range("M62").Value = -5726.8335
$ws.Range("N62").Value = -20750.875
$ws.Range("H65").Value = 13866.286
$ws.Range("I65").Value = 6350.8335
$ws.Range("J65").Value = 19502.875
$ws.Range("K65").Value = 31754.1675
$ws.Range("L65").Value = 97514.375
$ws.Range("M65").Value = -28634.1675
$ws.Range("N65").Value = -103754.375
$ws.Range("H86").Value = 5390.4
$ws.Range("I86").Value = 4984.6665
$ws.Range("K86").Value = 4984.6665
$ws.Range("M86").Value = -3861.6665
$ws.Range("H89").Value = 5390.4
$ws.Range("I89").Value = 4984.6665
$ws.Range("K89").Value = 24923.3325
$ws.Range("M89").Value = -19307.3325
$ws.Range("H122").Value = 3937.5908
$ws.Range("I122").Value = 4047.9375
$ws.Range("K122").Value = 12143.8125
$ws.Range("M122").Value = -9693.8125
$ws.Range("H134").Value = 6081.207
$ws.Range("I134").Value = 5470.5654
$ws.Range("K134").Value = 16411.6962
$ws.Range("M134").Value = -13876.6962

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 640.3333
$ws.Range("I14").Value = 640.3333
$ws.Range("K14").Value = 1920.9999
$ws.Range("M14").Value = -1747.9999
$ws.Range("H33").Value = 840.3333
$ws.Range("I33").Value = 1587.1428
$ws.Range("J33").Value = 186.875
$ws.Range("K33").Value = 9522.856800000001
$ws.Range("L33").Value = 1121.25
$ws.Range("M33").Value = -9239.856800000001
$ws.Range("N33").Value = -1687.25
$ws.Range("H68").Value = 1915.091
$ws.Range("I68").Value = 470.6
$ws.Range("J68").Value = 3118.8333
$ws.Range("K68").Value = 1411.8
$ws.Range("L68").Value = 9356.499899999999
$ws.Range("M68").Value = -600.8000000000002
$ws.Range("N68").Value = -10978.4999
$ws.Range("H71").Value = 1915.091
$ws.Range("I71").Value = 470.6
$ws.Range("J71").Value = 3118.8333
$ws.Range("K71").Value = 4235.400000000001
$ws.Range("L71").Value = 28069.4997
$ws.Range("M71").Value = -179.4000000000005
$ws.Range("N71").Value = -36181.4997
$ws.Range("H133").Value = 25834.785
$ws.Range("I133").Value = 18337.4
$ws.Range("K133").Value = 55012.2
$ws.Range("M133").Value = -49952.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4952.727
$ws.Range("I102").Value = 4430.385
$ws.Range("K102").Value = 4430.385
$ws.Range("M102").Value = -2808.385

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4078
$ws.Range("I16").Value = 4078
$ws.Range("K16").Value = 4078
$ws.Range("M16").Value = -3908
$ws.Range("H22").Value = 1429.7
$ws.Range("J22").Value = 1505.2778
$ws.Range("L22").Value = 1505.2778
$ws.Range("N22").Value = -2095.2778
$ws.Range("H27").Value = 1429.7
$ws.Range("J27").Value = 1505.2778
$ws.Range("L27").Value = 1505.2778
$ws.Range("N27").Value = -1719.2778
$ws.Range("H61").Value = 2488.6
$ws.Range("I61").Value = 2488.6
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2488.6
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2286.6
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 2488.6
$ws.Range("I113").Value = 2488.6
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2488.6
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -318.5999999999999
$ws.Range("N113").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 6250
$ws.Range("I9").Value = 2500
$ws.Range("K9").Value = 2500
$ws.Range("M9").Value = -2360
$ws.Range("H14").Value = 2507996.2
$ws.Range("J14").Value = 2995
$ws.Range("L14").Value = 2995
$ws.Range("N14").Value = -3331
$ws.Range("H126").Value = 1911.15
$ws.Range("I126").Value = 1479
$ws.Range("K126").Value = 4437
$ws.Range("M126").Value = -1967
